# Historia de Usuario 002 - arreglado creacion de nuevas hojas de excel
#
# 1) Append a block of new sample rows to "Java Books" (rows 20-32), where
#    row 20 and row 32 are blank spacer rows and rows 21-31 hold new data.
# 2) Add a second worksheet "Java Books 2" (right after "Java Books"),
#    with the same header row (No / Book Title / Author / Price).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# New rows appended at the bottom of "Java Books"
# ---------------------------------------------------------------------

# Row 20 - blank spacer row (kept present in the sheet, but with no cells)
$ws1.Rows.Item(20).OutlineLevel = 1
$ws1.Rows.Item(20).OutlineLevel = 0

# Row 21
$ws1.Range("B21").Value = 15.0
$ws1.Range("C21").Value = "Zzz"
$ws1.Range("D21").Value = "timmy"
$c = $ws1.Range("E21")
$c.Formula = "=TEXT(120,""000"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 22
$ws1.Range("A22").Value = 12.0
$ws1.Range("B22").Value = "fff"
$ws1.Range("C22").Value = "ccc"
$c = $ws1.Range("D22")
$c.Formula = "=TEXT(22,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 23
$ws1.Range("A23").Value = 11.0
$c = $ws1.Range("B23")
$c.Formula = "=TEXT(22,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)
$ws1.Range("C23").Value = "dd"
$c = $ws1.Range("D23")
$c.Formula = "=TEXT(33,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 24
$ws1.Range("A24").Value = 55.0
$ws1.Range("B24").Value = "ss"
$ws1.Range("C24").Value = "ss"
$c = $ws1.Range("D24")
$c.Formula = "=TEXT(44,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 25
$ws1.Range("A25").Value = 44.0
$ws1.Range("B25").Value = "uu"
$ws1.Range("C25").Value = "jj"
$c = $ws1.Range("D25")
$c.Formula = "=TEXT(77,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 26
$ws1.Range("A26").Value = 33.0
$c = $ws1.Range("B26")
$c.Formula = "=TEXT(33,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws1.Range("C26")
$c.Formula = "=TEXT(44,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws1.Range("D26")
$c.Formula = "=TEXT(55,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 27
$ws1.Range("A27").Value = 44.0
$c = $ws1.Range("B27")
$c.Formula = "=TEXT(55,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws1.Range("C27")
$c.Formula = "=TEXT(33,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws1.Range("D27")
$c.Formula = "=TEXT(66,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 28
$ws1.Range("A28").Value = 99.0
$ws1.Range("B28").Value = "nn"
$ws1.Range("C28").Value = "kk"
$c = $ws1.Range("D28")
$c.Formula = "=TEXT(0,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 29
$ws1.Range("A29").Value = 0.0
$ws1.Range("B29").Value = "nn"
$ws1.Range("C29").Value = "pp"
$ws1.Range("D29").Value = "bb"

# Row 30
$ws1.Range("A30").Value = 77.0
$ws1.Range("B30").Value = "nn"
$ws1.Range("C30").Value = "gg"
$c = $ws1.Range("D30")
$c.Formula = "=TEXT(99,""00"")"
$c.Copy()
$c.PasteSpecial(-4163)

# Row 31
$ws1.Range("A31").Value = 55.0
$ws1.Range("B31").Value = ".."
$ws1.Range("C31").Value = "jj"
$ws1.Range("D31").Value = "oo"

# Row 32 - trailing blank spacer row
$ws1.Rows.Item(32).OutlineLevel = 1
$ws1.Rows.Item(32).OutlineLevel = 0

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# New worksheet "Java Books 2", placed right after "Java Books",
# with the same header row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Java Books 2"
$ws2.Range("A1").Value = "No"
$ws2.Range("B1").Value = "Book Title"
$ws2.Range("C1").Value = "Author"
$ws2.Range("D1").Value = "Price"
